# Add a new worksheet '2025-11-26' with the weekly ranking data,
# placed after the last existing sheet, matching the style of prior weeks.
$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $srcSheet)
$newSheet.Name = "2025-11-26"

# Header row - copy formatting (bold, border, centered) from the previous week
$newSheet.Range("A1").Value = "rank"
$newSheet.Range("B1").Value = "title"
$newSheet.Range("C1").Value = "volume"
$newSheet.Range("D1").Value = "publisher"
$srcSheet.Range("A1:D1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$newSheet.Cells.Item(2, 1).Value = 1
$newSheet.Cells.Item(2, 2).Value = 'ブルーロック'
$newSheet.Cells.Item(2, 3).Value = 36
$newSheet.Cells.Item(3, 1).Value = 2
$newSheet.Cells.Item(3, 2).Value = '黄泉のツガイ'
$newSheet.Cells.Item(3, 3).Value = 11
$newSheet.Cells.Item(4, 1).Value = 3
$newSheet.Cells.Item(4, 2).Value = 'MIX'
$newSheet.Cells.Item(4, 3).Value = 24
$newSheet.Cells.Item(5, 1).Value = 4
$newSheet.Cells.Item(5, 2).Value = '転生賢者の異世界ライフ~第二の職業を得て、世界最強になりました~'
$newSheet.Cells.Item(5, 3).Value = 30
$newSheet.Cells.Item(6, 1).Value = 5
$newSheet.Cells.Item(6, 2).Value = '金色のガッシュ!! 2 Page'
$newSheet.Cells.Item(6, 3).Value = 35
$newSheet.Cells.Item(7, 1).Value = 6
$newSheet.Cells.Item(7, 2).Value = 'ゆびさきと恋々'
$newSheet.Cells.Item(7, 3).Value = 13
$newSheet.Cells.Item(8, 1).Value = 7
$newSheet.Cells.Item(8, 2).Value = 'ゆるキャン△'
$newSheet.Cells.Item(8, 3).Value = 18
$newSheet.Cells.Item(9, 1).Value = 8
$newSheet.Cells.Item(9, 2).Value = 'ワンパンマン'
$newSheet.Cells.Item(9, 3).Value = 35
$newSheet.Cells.Item(10, 1).Value = 9
$newSheet.Cells.Item(10, 2).Value = '出会って5秒でバトル'
$newSheet.Cells.Item(10, 3).Value = 29
$newSheet.Cells.Item(11, 1).Value = 10
$newSheet.Cells.Item(11, 2).Value = '彼女、お借りします'
$newSheet.Cells.Item(11, 3).Value = 43
$newSheet.Cells.Item(12, 1).Value = 11
$newSheet.Cells.Item(12, 2).Value = 'ガチアクタ'
$newSheet.Cells.Item(12, 3).Value = 17
$newSheet.Cells.Item(13, 1).Value = 12
$newSheet.Cells.Item(13, 2).Value = 'カッコウの許嫁'
$newSheet.Cells.Item(13, 3).Value = 30
$newSheet.Cells.Item(14, 1).Value = 13
$newSheet.Cells.Item(14, 2).Value = '裏バイト:逃亡禁止'
$newSheet.Cells.Item(14, 3).Value = 17
$newSheet.Cells.Item(15, 1).Value = 14
$newSheet.Cells.Item(15, 2).Value = 'みいちゃんと山田さん'
$newSheet.Cells.Item(15, 3).Value = 4
$newSheet.Cells.Item(16, 1).Value = 15
$newSheet.Cells.Item(16, 2).Value = '貸した魔力はで強制徴収~用済みとパーティー追放された俺は、可愛いサポート妖精と一緒に取り立てた魔力を運用して最強を目指す。~'
$newSheet.Cells.Item(16, 3).Value = 5
$newSheet.Cells.Item(17, 1).Value = 16
$newSheet.Cells.Item(17, 2).Value = '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐&『ざまぁ!』します!'
$newSheet.Cells.Item(17, 3).Value = 19
$newSheet.Cells.Item(18, 1).Value = 17
$newSheet.Cells.Item(18, 2).Value = '本好きの下剋上~司書になるためには手段を選んでいられません~第四部「貴族院の図書館を救いたい!11」'
$newSheet.Cells.Item(18, 3).Value = 11
$newSheet.Cells.Item(19, 1).Value = 18
$newSheet.Cells.Item(19, 2).Value = 'やんごとなき一族'
$newSheet.Cells.Item(19, 3).Value = 19
$newSheet.Cells.Item(20, 1).Value = 19
$newSheet.Cells.Item(20, 2).Value = '転生したら第七王子だったので、気ままに魔術を極めます'
$newSheet.Cells.Item(20, 3).Value = 21
$newSheet.Cells.Item(21, 1).Value = 20
$newSheet.Cells.Item(21, 2).Value = 'みいちゃんと山田さん'
$newSheet.Cells.Item(21, 3).Value = 3
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(21, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(21, 3).Value = 3
$newSheet.Cells.Item(22, 1).Value = 21
$newSheet.Cells.Item(22, 2).Value = 'みいちゃんと山田さん'
$newSheet.Cells.Item(22, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(22, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(22, 3).Value = 1
$newSheet.Cells.Item(23, 1).Value = 22
$newSheet.Cells.Item(23, 2).Value = 'みいちゃんと山田さん'
$newSheet.Cells.Item(23, 3).Value = 2
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(23, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(23, 3).Value = 2
$newSheet.Cells.Item(24, 1).Value = 23
$newSheet.Cells.Item(24, 2).Value = 'ダンス・ダンス・ダンスール'
$newSheet.Cells.Item(24, 3).Value = 31
$newSheet.Cells.Item(25, 1).Value = 24
$newSheet.Cells.Item(25, 2).Value = 'ちはやふる plus きみがため'
$newSheet.Cells.Item(25, 3).Value = 5
$newSheet.Cells.Item(26, 1).Value = 25
$newSheet.Cells.Item(26, 2).Value = '勘違いの工房主 英雄パーティの元雑用係が、実は戦闘以外がSSSランクだったというよくある話9'
$newSheet.Cells.Item(26, 3).Value = 9
$newSheet.Cells.Item(27, 1).Value = 26
$newSheet.Cells.Item(27, 2).Value = '虚構推理'
$newSheet.Cells.Item(27, 3).Value = 24
$newSheet.Cells.Item(28, 1).Value = 27
$newSheet.Cells.Item(28, 2).Value = '転生者は世間知らず ~特典スキルでスローライフ!?~ コミック版'
$newSheet.Cells.Item(28, 3).Value = 4
$newSheet.Cells.Item(29, 1).Value = 28
$newSheet.Cells.Item(29, 2).Value = '浪と損害のティティス'
$newSheet.Cells.Item(29, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(29, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(29, 3).Value = 1
$newSheet.Cells.Item(30, 1).Value = 29
$newSheet.Cells.Item(30, 2).Value = '暗殺後宮~暗殺女官・花鈴はゆったり生きたい~'
$newSheet.Cells.Item(30, 3).Value = 9
$newSheet.Cells.Item(31, 1).Value = 30
$newSheet.Cells.Item(31, 2).Value = '紫雲寺家の子供たち'
$newSheet.Cells.Item(31, 3).Value = 8
$newSheet.Cells.Item(32, 1).Value = 31
$newSheet.Cells.Item(32, 2).Value = 'DEAR BOYS ACT4'
$newSheet.Cells.Item(32, 3).Value = 21
$newSheet.Cells.Item(33, 1).Value = 32
$newSheet.Cells.Item(33, 2).Value = '戦隊大失格'
$newSheet.Cells.Item(33, 3).Value = 20
$newSheet.Cells.Item(34, 1).Value = 33
$newSheet.Cells.Item(34, 2).Value = '杖と剣のウィストリア'
$newSheet.Cells.Item(34, 3).Value = 14
$newSheet.Cells.Item(35, 1).Value = 34
$newSheet.Cells.Item(35, 2).Value = 'ハーレム・メイカー~ゲームのヒロインたちの攻略対象が俺ってマジか?~'
$newSheet.Cells.Item(35, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(35, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(35, 3).Value = 1
$newSheet.Cells.Item(36, 1).Value = 35
$newSheet.Cells.Item(36, 2).Value = 'アイツ'
$newSheet.Cells.Item(36, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(36, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(36, 3).Value = 1
$newSheet.Cells.Item(37, 1).Value = 36
$newSheet.Cells.Item(37, 2).Value = 'BLUE GIANT MOMENTUM'
$newSheet.Cells.Item(37, 3).Value = 6
$newSheet.Cells.Item(38, 1).Value = 37
$newSheet.Cells.Item(38, 2).Value = 'KING GOLF'
$newSheet.Cells.Item(38, 3).Value = 43
$newSheet.Cells.Item(39, 1).Value = 38
$newSheet.Cells.Item(39, 2).Value = 'SPY×FAMILY'
$newSheet.Cells.Item(39, 3).Value = 16
$newSheet.Cells.Item(40, 1).Value = 39
$newSheet.Cells.Item(40, 2).Value = '魔導具師ダリヤはうつむかない ~王立高等学院編~'
$newSheet.Cells.Item(40, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(40, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(40, 3).Value = 1
$newSheet.Cells.Item(41, 1).Value = 40
$newSheet.Cells.Item(41, 2).Value = '復讐の輪廻2'
$newSheet.Cells.Item(41, 3).Value = 2
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(41, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(41, 3).Value = 2
$newSheet.Cells.Item(42, 1).Value = 41
$newSheet.Cells.Item(42, 2).Value = 'ユウリ~彼女が復讐する理由~2'
$newSheet.Cells.Item(42, 3).Value = 2
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(42, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(42, 3).Value = 2
$newSheet.Cells.Item(43, 1).Value = 42
$newSheet.Cells.Item(43, 2).Value = 'kitty,kitty,kitty! -ケダモノアラシ-'
$newSheet.Cells.Item(43, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(43, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(43, 3).Value = 1
$newSheet.Cells.Item(44, 1).Value = 43
$newSheet.Cells.Item(44, 2).Value = '未知と宝物ざっくざくの迷宮大配信! ~ハズレスキルすらない凡人、見る人から見れば普通に非凡でした~ コミック版'
$newSheet.Cells.Item(44, 3).Value = 4
$newSheet.Cells.Item(45, 1).Value = 44
$newSheet.Cells.Item(45, 2).Value = '四姉妹は夜をおまちかね'
$newSheet.Cells.Item(45, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(45, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(45, 3).Value = 1
$newSheet.Cells.Item(46, 1).Value = 45
$newSheet.Cells.Item(46, 2).Value = '浪と損害のティティス'
$newSheet.Cells.Item(46, 3).Value = 2
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(46, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(46, 3).Value = 2
$newSheet.Cells.Item(47, 1).Value = 46
$newSheet.Cells.Item(47, 2).Value = '花束のような復讐を~怪物たちの愛を知れ~'
$newSheet.Cells.Item(47, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(47, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(47, 3).Value = 1
$newSheet.Cells.Item(48, 1).Value = 47
$newSheet.Cells.Item(48, 2).Value = '無能の中の無能王子 スキルを授かりましたが、周りの女性はとかです コミック版'
$newSheet.Cells.Item(48, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(48, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(48, 3).Value = 1
$newSheet.Cells.Item(49, 1).Value = 48
$newSheet.Cells.Item(49, 2).Value = 'わたしの地味セン王子'
$newSheet.Cells.Item(49, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(49, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(49, 3).Value = 1
$newSheet.Cells.Item(50, 1).Value = 49
$newSheet.Cells.Item(50, 2).Value = '硝子白書'
$newSheet.Cells.Item(50, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(50, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(50, 3).Value = 1
$newSheet.Cells.Item(51, 1).Value = 50
$newSheet.Cells.Item(51, 2).Value = '瞳いっぱいの涙'
$newSheet.Cells.Item(51, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(51, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(51, 3).Value = 1
$newSheet.Cells.Item(52, 1).Value = 51
$newSheet.Cells.Item(52, 2).Value = '菜子の色'
$newSheet.Cells.Item(52, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(52, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(52, 3).Value = 1
$newSheet.Cells.Item(53, 1).Value = 52
$newSheet.Cells.Item(53, 2).Value = '君だけに輝く'
$newSheet.Cells.Item(53, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(53, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(53, 3).Value = 1
$newSheet.Cells.Item(54, 1).Value = 53
$newSheet.Cells.Item(54, 2).Value = '賭ケグルイ'
$newSheet.Cells.Item(54, 3).Value = 20
$newSheet.Cells.Item(55, 1).Value = 54
$newSheet.Cells.Item(55, 2).Value = '火の神さまの掃除人ですが、いつの間にか花嫁として溺愛されています'
$newSheet.Cells.Item(55, 3).Value = 7
$newSheet.Cells.Item(56, 1).Value = 55
$newSheet.Cells.Item(56, 2).Value = '野生のラスボスが現れた! 黒翼の覇王11'
$newSheet.Cells.Item(56, 3).Value = 11
$newSheet.Cells.Item(57, 1).Value = 56
$newSheet.Cells.Item(57, 2).Value = '転生したら皇帝でした~生まれながらの皇帝はこの先生き残れるか~@COMIC'
$newSheet.Cells.Item(57, 3).Value = 6
$newSheet.Cells.Item(58, 1).Value = 57
$newSheet.Cells.Item(58, 2).Value = '時々ボソッとロシア語でデレる隣のアーリャさん'
$newSheet.Cells.Item(58, 3).Value = 8
$newSheet.Cells.Item(59, 1).Value = 58
$newSheet.Cells.Item(59, 2).Value = '灰仭巫覡'
$newSheet.Cells.Item(59, 3).Value = 6
$newSheet.Cells.Item(60, 1).Value = 59
$newSheet.Cells.Item(60, 2).Value = '蒼く染めろ'
$newSheet.Cells.Item(60, 3).Value = 18
$newSheet.Cells.Item(61, 1).Value = 60
$newSheet.Cells.Item(61, 2).Value = '憂国のモリアーティ'
$newSheet.Cells.Item(61, 3).Value = 21
$newSheet.Cells.Item(62, 1).Value = 61
$newSheet.Cells.Item(62, 2).Value = '服飾師ルチアはあきらめない ~今日から始める幸服計画~'
$newSheet.Cells.Item(62, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(62, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(62, 3).Value = 1
$newSheet.Cells.Item(63, 1).Value = 62
$newSheet.Cells.Item(63, 2).Value = '独身偽装~私の彼氏は既婚者でした~2'
$newSheet.Cells.Item(63, 3).Value = 2
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(63, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(63, 3).Value = 2
$newSheet.Cells.Item(64, 1).Value = 63
$newSheet.Cells.Item(64, 2).Value = '俺たちの善と恋について'
$newSheet.Cells.Item(64, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(64, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(64, 3).Value = 1
$newSheet.Cells.Item(65, 1).Value = 64
$newSheet.Cells.Item(65, 2).Value = '浪と損害のティティス'
$newSheet.Cells.Item(65, 3).Value = 3
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(65, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(65, 3).Value = 3
$newSheet.Cells.Item(66, 1).Value = 65
$newSheet.Cells.Item(66, 2).Value = '声が聞きたい小平くん'
$newSheet.Cells.Item(66, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(66, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(66, 3).Value = 1
$newSheet.Cells.Item(67, 1).Value = 66
$newSheet.Cells.Item(67, 2).Value = '一年後、生贄になる君と偽りの恋をする'
$newSheet.Cells.Item(67, 3).Value = 2
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(67, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(67, 3).Value = 2
$newSheet.Cells.Item(68, 1).Value = 67
$newSheet.Cells.Item(68, 2).Value = 'いつわりの花嫁 ~旦那さま、今宵お命頂戴します~'
$newSheet.Cells.Item(68, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(68, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(68, 3).Value = 1
$newSheet.Cells.Item(69, 1).Value = 68
$newSheet.Cells.Item(69, 2).Value = 'わたしの地味セン王子'
$newSheet.Cells.Item(69, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(69, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(69, 3).Value = 1
$newSheet.Cells.Item(70, 1).Value = 69
$newSheet.Cells.Item(70, 2).Value = 'Friends -制服イレブン-'
$newSheet.Cells.Item(70, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(70, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(70, 3).Value = 1
$newSheet.Cells.Item(71, 1).Value = 70
$newSheet.Cells.Item(71, 2).Value = 'Believe'
$newSheet.Cells.Item(71, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(71, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(71, 3).Value = 1
$newSheet.Cells.Item(72, 1).Value = 71
$newSheet.Cells.Item(72, 2).Value = 'てぃーんず -制服の林檎たち-'
$newSheet.Cells.Item(72, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(72, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(72, 3).Value = 1
$newSheet.Cells.Item(73, 1).Value = 72
$newSheet.Cells.Item(73, 2).Value = '装備製作系チートで異世界を自由に生きていきます6'
$newSheet.Cells.Item(73, 3).Value = 6
$newSheet.Cells.Item(74, 1).Value = 73
$newSheet.Cells.Item(74, 2).Value = '転生したらスライムだった件'
$newSheet.Cells.Item(74, 3).Value = 30
$newSheet.Cells.Item(75, 1).Value = 74
$newSheet.Cells.Item(75, 2).Value = 'ババンババンバンバンパイア'
$newSheet.Cells.Item(75, 3).Value = 12
$newSheet.Cells.Item(76, 1).Value = 75
$newSheet.Cells.Item(76, 2).Value = 'めっちゃ召喚された件 THE COMIC'
$newSheet.Cells.Item(76, 3).Value = 11
$newSheet.Cells.Item(77, 1).Value = 76
$newSheet.Cells.Item(77, 2).Value = 'SAKAMOTO DAYS'
$newSheet.Cells.Item(77, 3).Value = 24
$newSheet.Cells.Item(78, 1).Value = 77
$newSheet.Cells.Item(78, 2).Value = '無職転生 ~異世界行ったら本気だす~ 失意の魔術師編'
$newSheet.Cells.Item(78, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(78, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(78, 3).Value = 1
$newSheet.Cells.Item(79, 1).Value = 78
$newSheet.Cells.Item(79, 2).Value = 'フリースキルで最強冒険者 ~ペットも無双で異世界生活が楽しすぎる~'
$newSheet.Cells.Item(79, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(79, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(79, 3).Value = 1
$newSheet.Cells.Item(80, 1).Value = 79
$newSheet.Cells.Item(80, 2).Value = 'お茶屋さんは賢者見習い'
$newSheet.Cells.Item(80, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(80, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(80, 3).Value = 1
$newSheet.Cells.Item(81, 1).Value = 80
$newSheet.Cells.Item(81, 2).Value = '婚約破棄されたらエリート御曹司の義弟に娶られました1'
$newSheet.Cells.Item(81, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(81, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(81, 3).Value = 1
$newSheet.Cells.Item(82, 1).Value = 81
$newSheet.Cells.Item(82, 2).Value = 'シンデレラ・コンプレックス 1話 始まりの教室'
$newSheet.Cells.Item(82, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(82, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(82, 3).Value = 1
$newSheet.Cells.Item(83, 1).Value = 82
$newSheet.Cells.Item(83, 2).Value = 'これが恋だと知っている'
$newSheet.Cells.Item(83, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(83, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(83, 3).Value = 1
$newSheet.Cells.Item(84, 1).Value = 83
$newSheet.Cells.Item(84, 2).Value = '現実世界にダンジョン現る! ~アラサーフリーターは元聖女のスケルトンと一緒に成り上がります!~ コミック版'
$newSheet.Cells.Item(84, 3).Value = 3
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(84, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(84, 3).Value = 3
$newSheet.Cells.Item(85, 1).Value = 84
$newSheet.Cells.Item(85, 2).Value = 'RED&BLUE 第7話'
$newSheet.Cells.Item(85, 3).Value = 7
$newSheet.Cells.Item(86, 1).Value = 85
$newSheet.Cells.Item(86, 2).Value = '経験人数が見えるメガネ'
$newSheet.Cells.Item(86, 3).Value = 4
$newSheet.Cells.Item(87, 1).Value = 86
$newSheet.Cells.Item(87, 2).Value = 'オンタマ!?ミラクルフラッピン!'
$newSheet.Cells.Item(87, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(87, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(87, 3).Value = 1
$newSheet.Cells.Item(88, 1).Value = 87
$newSheet.Cells.Item(88, 2).Value = '恋はカーテンコールのあとで。'
$newSheet.Cells.Item(88, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(88, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(88, 3).Value = 1
$newSheet.Cells.Item(89, 1).Value = 88
$newSheet.Cells.Item(89, 2).Value = 'いつわり婚'
$newSheet.Cells.Item(89, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(89, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(89, 3).Value = 1
$newSheet.Cells.Item(90, 1).Value = 89
$newSheet.Cells.Item(90, 2).Value = '召喚聖女は魔王様の膝の上~聖女は猫になりまして~'
$newSheet.Cells.Item(90, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(90, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(90, 3).Value = 1
$newSheet.Cells.Item(91, 1).Value = 90
$newSheet.Cells.Item(91, 2).Value = 'レベルアップデイズ~幼馴染の攻略サポート~'
$newSheet.Cells.Item(91, 3).Value = 2
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(91, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(91, 3).Value = 2
$newSheet.Cells.Item(92, 1).Value = 91
$newSheet.Cells.Item(92, 2).Value = 'すばらしき新世界(フルカラー)'
$newSheet.Cells.Item(92, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(92, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(92, 3).Value = 1
$newSheet.Cells.Item(93, 1).Value = 92
$newSheet.Cells.Item(93, 2).Value = 'かぐや姫の孫'
$newSheet.Cells.Item(93, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(93, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(93, 3).Value = 1
$newSheet.Cells.Item(94, 1).Value = 93
$newSheet.Cells.Item(94, 2).Value = '現実世界にダンジョン現る! ~アラサーフリーターは元聖女のスケルトンと一緒に成り上がります!~ コミック版'
$newSheet.Cells.Item(94, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(94, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(94, 3).Value = 1
$newSheet.Cells.Item(95, 1).Value = 94
$newSheet.Cells.Item(95, 2).Value = '推しの敵になったので@COMIC 第1話'
$newSheet.Cells.Item(95, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(95, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(95, 3).Value = 1
$newSheet.Cells.Item(96, 1).Value = 95
$newSheet.Cells.Item(96, 2).Value = '芦屋山手 お道具迎賓館(コミック) 1話'
$newSheet.Cells.Item(96, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(96, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(96, 3).Value = 1
$newSheet.Cells.Item(97, 1).Value = 96
$newSheet.Cells.Item(97, 2).Value = '入れ替わったら、オレ様彼氏とエッチする運命でした!'
$newSheet.Cells.Item(97, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(97, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(97, 3).Value = 1
$newSheet.Cells.Item(98, 1).Value = 97
$newSheet.Cells.Item(98, 2).Value = 'ウィズ -幸せのある場所-'
$newSheet.Cells.Item(98, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(98, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(98, 3).Value = 1
$newSheet.Cells.Item(99, 1).Value = 98
$newSheet.Cells.Item(99, 2).Value = '天使のオシャベリ'
$newSheet.Cells.Item(99, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(99, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(99, 3).Value = 1
$newSheet.Cells.Item(100, 1).Value = 99
$newSheet.Cells.Item(100, 2).Value = 'Kissの奇跡'
$newSheet.Cells.Item(100, 3).Value = 1
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(100, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(100, 3).Value = 1
$newSheet.Cells.Item(101, 1).Value = 100
$newSheet.Cells.Item(101, 2).Value = '君はぼくのヒーローさ'
$newSheet.Cells.Item(101, 3).Value = 3
$srcSheet.Range("C29").Copy()
$newSheet.Cells.Item(101, 3).PasteSpecial(-4122)
$newSheet.Cells.Item(101, 3).Value = 3

# Restore the originally active sheet/selection (adding a sheet shouldn't
# change which tab the workbook opens to)
$wb.Worksheets.Item(1).Activate()

Write-Output ("Created sheet: " + $newSheet.Name)
